$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.232.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.445.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.892.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.089.93"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.443.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +12.94%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "614.81"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.571.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.17%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.77%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.89"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.17"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0533"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.601"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0232"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.39%  "
